# Commit: "committing QueueTest.java, RegisterTest.java,Utils.java,excel sheet,DataProviderClass.java"
# Adds four new worksheets (registerpage, logindata, validcode, invalidcode) with
# register/login test data, mailto hyperlinks on the "@"-style credential strings,
# and a couple of rich-text (multi-font) cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the four new worksheets, in order, at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$registerpage = $wb.Worksheets.Add($null, $lastSheet)
$registerpage.Name = "registerpage"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$logindata = $wb.Worksheets.Add($null, $lastSheet)
$logindata.Name = "logindata"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$validcode = $wb.Worksheets.Add($null, $lastSheet)
$validcode.Name = "validcode"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$invalidcode = $wb.Worksheets.Add($null, $lastSheet)
$invalidcode.Name = "invalidcode"

Write-Host "Sheets added. Count = " $wb.Worksheets.Count

# ---------------------------------------------------------------------------
# 2. registerpage
# ---------------------------------------------------------------------------
$ws = $registerpage

$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "passwordconfirmation"
$ws.Range("D1").Value = "errormsg"

$ws.Range("A2").Value = "Lillyy_1@991"
$ws.Range("B2").Value = "testPassword@1"
$ws.Range("C2").Value = "testPassword@1 "
$ws.Range("D2").Value = "username already exists"

$ws.Range("A3").Value = "Raha_a@123"
$ws.Range("B3").Value = "ah_1@4"
$ws.Range("C3").Value = "ah_1@4"
# D3: rich text - "Password should be atleast " + "8" (green) + " characters " (black)
$ws.Range("D3").Value = "Password should be atleast 8 characters "

$ws.Range("A4").Value = "&&**&&**"
$ws.Range("B4").Value = "aha_1@456"
$ws.Range("C4").Value = "aha_1@456"
$ws.Range("D4").Value = "Please enter a valid username "

$ws.Range("A5").Value = "Raha_a@123"
$ws.Range("B5").Value = 12345678
$ws.Range("C5").Value = 12345678
$ws.Range("D5").Value = "Password cannot be only numeric "

$ws.Range("A6").Value = "Raha_a@123"
$ws.Range("B6").Value = "aha_1@456"
$ws.Range("C6").Value = "aha_1@4567 "
$ws.Range("D6").Value = [char]0x2019
$ws.Range("D6").Value = "password_mismatch:The two password fields didn" + [char]0x2019 + "t match."

# Borders A2:D6, thin box around every cell (header row A1:D1 stays unformatted)
$ws.Range("A2:D6").Borders.LineStyle = 1

# Row height 15.5 for rows 2-5 (ht="15.5" in target XML; row 6 stays default)
$ws.Range("A2:D5").RowHeight = 15.5

# Vertical alignment = top for rows 2-6
$ws.Range("A2:D6").VerticalAlignment = -4160

# Base font for column D (errormsg) + A4: Consolas 12 black, except D6 stays default font
$ws.Range("D2:D5").Font.Name = "Consolas"
$ws.Range("D2:D5").Font.Size = 12
$ws.Range("D2:D5").Font.Color = 0
$ws.Range("A4").Font.Name = "Consolas"
$ws.Range("A4").Font.Size = 12
$ws.Range("A4").Font.Color = 0

# B5/C5 numeric cells: Consolas 12 green
$ws.Range("B5:C5").Font.Name = "Consolas"
$ws.Range("B5:C5").Font.Size = 12
$ws.Range("B5:C5").Font.Color = 32768

# Column widths (approximate best-fit)
$ws.Columns.Item(1).ColumnWidth = 11.8
$ws.Columns.Item(2).ColumnWidth = 14.8
$ws.Columns.Item(3).ColumnWidth = 19.8
$ws.Columns.Item(4).ColumnWidth = 51.1

# Hyperlinks (Excel auto-mailto style for "@"-shaped credential strings)
$null = $ws.Hyperlinks.Add($ws.Range("A2"), "mailto:Lillyy_1@991")
$null = $ws.Hyperlinks.Add($ws.Range("B2"), "mailto:testPassword@1")
$null = $ws.Hyperlinks.Add($ws.Range("C2"), "mailto:testPassword@1")
$null = $ws.Hyperlinks.Add($ws.Range("A3"), "mailto:Raha_a@123")
$null = $ws.Hyperlinks.Add($ws.Range("B3"), "mailto:ah_1@4")
$null = $ws.Hyperlinks.Add($ws.Range("C3"), "mailto:ah_1@4")
$null = $ws.Hyperlinks.Add($ws.Range("B4"), "mailto:aha_1@456")
$null = $ws.Hyperlinks.Add($ws.Range("C4"), "mailto:aha_1@456")
$null = $ws.Hyperlinks.Add($ws.Range("A5"), "mailto:Raha_a@123")
$null = $ws.Hyperlinks.Add($ws.Range("A6"), "mailto:Raha_a@123")
$null = $ws.Hyperlinks.Add($ws.Range("B6"), "mailto:aha_1@456")
$null = $ws.Hyperlinks.Add($ws.Range("C6"), "mailto:aha_1@4567")

# Rich text runs -------------------------------------------------------
# D3: "Password should be atleast " + "8"(green Consolas12) + " characters "(black Consolas12)
$d3 = $ws.Range("D3")
$run2start = ("Password should be atleast ").Length + 1
$run2len = ("8").Length
$run3start = $run2start + $run2len
$run3len = (" characters ").Length
$c = $d3.Characters($run2start, $run2len)
$c.Font.Name = "Consolas"
$c.Font.Size = 12
$c.Font.Color = 32768
$c = $d3.Characters($run3start, $run3len)
$c.Font.Name = "Consolas"
$c.Font.Size = 12
$c.Font.Color = 0

# A4: "&&" (default font) + "**&&**" (teal Consolas12)
$a4 = $ws.Range("A4")
$c = $a4.Characters(3, 6)
$c.Font.Name = "Consolas"
$c.Font.Size = 12
$c.Font.Color = 8421376

$ws.Cells.Select()

# ---------------------------------------------------------------------------
# 3. logindata
# ---------------------------------------------------------------------------
$ws = $logindata
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "Lillyy_1@991"
$ws.Range("B2").Value = "testPassword@1"

$ws.Range("A2:B2").RowHeight = 15.5
$ws.Range("A2:B2").Font.Name = "Consolas"
$ws.Range("A2:B2").Font.Size = 12
$ws.Range("A2:B2").Font.Color = 16711722

$ws.Columns.Item(1).ColumnWidth = 15.3
$ws.Columns.Item(2).ColumnWidth = 17.7

$ws.Cells.Select()

